$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '29.495.59'
$ws.Cells.Item(2, 5).Value = '  +2.15%  '

$ws.Cells.Item(3, 4).Value = '1.918.58'
$ws.Cells.Item(3, 5).Value = '  +1.66%  '

$ws.Cells.Item(4, 5).Value = '  -0.20%  '

$ws.Cells.Item(5, 4).Value = '325.94'
$ws.Cells.Item(5, 5).Value = '  -1.78%  '

$ws.Cells.Item(6, 4).Value = '1.000'
$ws.Cells.Item(6, 5).Value = '  -0.18%  '

$ws.Cells.Item(7, 4).Value = '0.4741'

$ws.Cells.Item(8, 4).Value = '0.4095'
$ws.Cells.Item(8, 5).Value = '  -0.40%  '

$ws.Cells.Item(9, 5).Value = '  +0.64%  '

$ws.Cells.Item(10, 4).Value = '0.08052'
$ws.Cells.Item(10, 5).Value = '  +0.98%  '

$ws.Cells.Item(11, 5).Value = '  +1.83%  '

$ws.Cells.Item(12, 4).Value = '22.50'

$ws.Cells.Item(13, 4).Value = '1.928.16'
$ws.Cells.Item(13, 5).Value = '  +1.74%  '

$ws.Cells.Item(14, 4).Value = '5.933'
$ws.Cells.Item(14, 5).Value = '  +0.36%  '

$ws.Cells.Item(15, 4).Value = '7.159'
$ws.Cells.Item(15, 5).Value = '  +1.38%  '

$ws.Cells.Item(16, 4).Value = '89.60'
$ws.Cells.Item(16, 5).Value = '  +0.35%  '

$ws.Cells.Item(17, 5).Value = '  -0.22%  '

$ws.Cells.Item(18, 2).Value = 'ShibaInu'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Cells.Item(18, 4).Value = '0.00001033'
$ws.Cells.Item(18, 5).Value = '  +0.46%  '

$ws.Cells.Item(19, 2).Value = 'TRON'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Cells.Item(19, 4).Value = '0.06598'
$ws.Cells.Item(19, 5).Value = '  +0.44%  '

$ws.Cells.Item(20, 4).Value = '17.80'
$ws.Cells.Item(20, 5).Value = '  +1.84%  '

$ws.Cells.Item(21, 4).Value = '1.0000'
$ws.Cells.Item(21, 5).Value = '  -0.19%  '

$ws.Cells.Item(22, 4).Value = '29.503.66'
$ws.Cells.Item(22, 5).Value = '  +1.97%  '

$ws.Cells.Item(23, 4).Value = '5.549'
$ws.Cells.Item(23, 5).Value = '  +3.01%  '

$ws.Cells.Item(24, 5).Value = '  +2.14%  '

$ws.Cells.Item(25, 4).Value = '2.208'
$ws.Cells.Item(25, 5).Value = '  -0.36%  '

$ws.Cells.Item(26, 4).Value = '2.132.96'
$ws.Cells.Item(26, 5).Value = '  +0.62%  '

$ws.Cells.Item(27, 4).Value = '154.73'

$ws.Cells.Item(28, 4).Value = '19.84'
$ws.Cells.Item(28, 5).Value = '  +0.83%  '

$ws.Cells.Item(29, 4).Value = '6.070'
$ws.Cells.Item(29, 5).Value = '  +12.00%  '

$ws.Cells.Item(30, 4).Value = '2.130'
$ws.Cells.Item(30, 5).Value = '  +0.52%  '

$ws.Cells.Item(31, 4).Value = '117.90'
$ws.Cells.Item(31, 5).Value = '  +0.01%  '

$ws.Cells.Item(32, 5).Value = '  +8.69%  '

$ws.Cells.Item(33, 4).Value = '0.09537'
$ws.Cells.Item(33, 5).Value = '  +1.88%  '

$ws.Cells.Item(34, 5).Value = '  +1.39%  '

$ws.Cells.Item(35, 5).Value = '  -1.21%  '

$ws.Cells.Item(36, 4).Value = '5.410'
$ws.Cells.Item(36, 5).Value = '  +2.49%  '

$ws.Cells.Item(37, 4).Value = '0.06121'
$ws.Cells.Item(37, 5).Value = '  +0.97%  '

$ws.Cells.Item(38, 5).Value = '  +1.20%  '

$ws.Cells.Item(39, 4).Value = '8.328'
$ws.Cells.Item(39, 5).Value = '  +0.38%  '

$ws.Cells.Item(40, 4).Value = '1.172'
$ws.Cells.Item(40, 5).Value = '  -0.35%  '

$ws.Cells.Item(41, 4).Value = '0.5897'
$ws.Cells.Item(41, 5).Value = '  +2.15%  '

$ws.Cells.Item(42, 4).Value = '2.560'
$ws.Cells.Item(42, 5).Value = '  +12.14%  '

$ws.Cells.Item(43, 4).Value = '0.1846'
$ws.Cells.Item(43, 5).Value = '  +1.64%  '

$ws.Cells.Item(44, 5).Value = '  +0.19%  '

$ws.Cells.Item(45, 4).Value = '0.08020'
$ws.Cells.Item(45, 5).Value = '  +14.37%  '

$ws.Cells.Item(46, 4).Value = '1.286'
$ws.Cells.Item(46, 5).Value = '  +1.67%  '

$ws.Cells.Item(47, 4).Value = '0.5564'
$ws.Cells.Item(47, 5).Value = '  +1.40%  '

$ws.Cells.Item(48, 4).Value = '12.12'
$ws.Cells.Item(48, 5).Value = '  +1.12%  '

$ws.Cells.Item(49, 4).Value = '1.939'
$ws.Cells.Item(49, 5).Value = '  +1.76%  '

$ws.Cells.Item(50, 4).Value = '113.13'
$ws.Cells.Item(50, 5).Value = '  +2.06%  '

$ws.Cells.Item(51, 4).Value = '45.18'
$ws.Cells.Item(51, 5).Value = '  +0.72%  '
